# Generate Report for Archive
# Update localization status for the two files that have moved from
# "Ready for handoff" to "In Translation": 19fe273b-...md and 24a188bc-...md

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns B (zh-cn) and C (de-de) show the status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("B4").Value = $newStatus
$overview.Range("C4").Value = $newStatus

# --- zh-cn sheet: column C is "Status" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

# --- de-de sheet: column C is "Status" ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus
